# Re-sort the language/value table in descending order by value (column B)
# and drop the two rows (Swedish, Uzbek) that no longer belong in the
# finished "imf gdp nominal multiyear" table, shrinking the used range
# from A1:B23 down to A1:B21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final ordered (language, value) pairs for rows 2..21, sorted descending by value.
$data = @(
    @("English", 21.71078478681658),
    @("Chinese", 18.03949587161287),
    @("Spanish", 6.780351875257941),
    @("German", 4.373171315329072),
    @("Arabic", 4.291748648722796),
    @("Japanese", 4.167346530704183),
    @("Russian", 3.23865925344384),
    @("Malay-Indonesian", 3.087401045354373),
    @("Portuguese", 2.90341152005385),
    @("French", 2.636550167086114),
    @("Italian", 2.01410968232425),
    @("Turkish", 1.854264366746624),
    @("Korean", 1.709734241349971),
    @("Dutch", 1.236803880441917),
    @("Persian", 1.051995398093134),
    @("Thai", 0.9915127621281781),
    @("Polish", 0.9539285658512771),
    @("Urdu", 0.8859759322005877),
    @("Vietnamese", 0.7215643340029804),
    @("Bengali", 0.7085157210901017)
)

# Remove the two trailing rows (22 and 23) that fall out of the final table.
$ws.Rows(23).Delete()
$ws.Rows(22).Delete()

# Write the re-sorted data into rows 2..21.
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
